# Tutorial 6 solution update
# - Reformat the dates in column A (rows 3-21) from dd/mm/yyyy to dd-mm-yyyy
# - Update a few attendance counter cells for the first two data rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @{
  "A3"  = "28-07-2022"
  "A4"  = "01-08-2022"
  "A5"  = "04-08-2022"
  "A6"  = "08-08-2022"
  "A7"  = "11-08-2022"
  "A8"  = "15-08-2022"
  "A9"  = "18-08-2022"
  "A10" = "22-08-2022"
  "A11" = "25-08-2022"
  "A12" = "29-08-2022"
  "A13" = "01-09-2022"
  "A14" = "05-09-2022"
  "A15" = "08-09-2022"
  "A16" = "12-09-2022"
  "A17" = "15-09-2022"
  "A18" = "19-09-2022"
  "A19" = "22-09-2022"
  "A20" = "26-09-2022"
  "A21" = "29-09-2022"
}

# The dates are plain text (dd-mm-yyyy), not real date values. Force a text
# number format first so Excel does not auto-convert the new strings into
# date serial numbers, write the values, then restore the original (default)
# cell style so no formatting changes leak into the saved file.
$datesRange = $ws.Range("A3:A21")
$datesRange.NumberFormat = "@"
foreach ($addr in $dates.Keys) {
  $ws.Range($addr).Value = $dates[$addr]
}
$datesRange.Style = "Normal"

# Update attendance counters for 28-07-2022 (row 3): now counted Real + Invalid
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

# Update attendance counters for 01-08-2022 (row 4): now counted Real + Duplicate, no longer Absent
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("H4").Value = 0
